$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 62503624
$ws.Range("I51").Value = 83336080
$ws.Range("K51").Value = 83336080
$ws.Range("M51").Value = -83335596
$ws.Range("H58").Value = 1306.625
$ws.Range("J58").Value = 495
$ws.Range("L58").Value = 1485
$ws.Range("N58").Value = -1785
$ws.Range("H99").Value = 358.45834
$ws.Range("I99").Value = 258.2381
$ws.Range("J99").Value = 1060
$ws.Range("K99").Value = 774.7142999999999
$ws.Range("L99").Value = 3180
$ws.Range("M99").Value = 723.2857000000001
$ws.Range("N99").Value = -6176
$ws.Range("H100").Value = 749.8
$ws.Range("I100").Value = 749.8
$ws.Range("K100").Value = 749.8
$ws.Range("M100").Value = -208.8
$ws.Range("H137").Value = 3389.946
$ws.Range("I137").Value = 2633.276
$ws.Range("K137").Value = 7899.828
$ws.Range("M137").Value = -5349.828

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1604.3871
$ws.Range("I32").Value = 1603
$ws.Range("K32").Value = 1603
$ws.Range("M32").Value = -1316
$ws.Range("H61").Value = 2476.111
$ws.Range("I61").Value = 2040.7142
$ws.Range("K61").Value = 2040.7142
$ws.Range("M61").Value = -1828.7142
$ws.Range("H136").Value = 2476.111
$ws.Range("I136").Value = 2040.7142
$ws.Range("K136").Value = 6122.142599999999
$ws.Range("M136").Value = -3572.142599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 718.6842
$ws.Range("I22").Value = 773.4706
$ws.Range("K22").Value = 773.4706
$ws.Range("M22").Value = -600.4706
$ws.Range("H107").Value = 2199121.2
$ws.Range("I107").Value = 2850195.5
$ws.Range("K107").Value = 2850195.5
$ws.Range("M107").Value = -2848275.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2235676.8
$ws.Range("I31").Value = 3346.7727
$ws.Range("J31").Value = 3680125.5
$ws.Range("K31").Value = 3346.7727
$ws.Range("L31").Value = 3680125.5
$ws.Range("M31").Value = -3051.7727
$ws.Range("N31").Value = -3680715.5
$ws.Range("H34").Value = 2235676.8
$ws.Range("I34").Value = 3346.7727
$ws.Range("J34").Value = 3680125.5
$ws.Range("K34").Value = 3346.7727
$ws.Range("L34").Value = 3680125.5
$ws.Range("M34").Value = -3144.7727
$ws.Range("N34").Value = -3680529.5
$ws.Range("H58").Value = 2092.0322
$ws.Range("I58").Value = 1495.5555
$ws.Range("J58").Value = 2917.923
$ws.Range("K58").Value = 1495.5555
$ws.Range("L58").Value = 2917.923
$ws.Range("M58").Value = -1292.5555
$ws.Range("N58").Value = -3323.923
$ws.Range("H105").Value = 3841.2856
$ws.Range("I105").Value = 2796.3333
$ws.Range("J105").Value = 4625
$ws.Range("K105").Value = 2796.3333
$ws.Range("L105").Value = 4625
$ws.Range("M105").Value = -1049.3333
$ws.Range("N105").Value = -8119
$ws.Range("H122").Value = 1120.6522
$ws.Range("I122").Value = 1152.75
$ws.Range("J122").Value = 906.6667
$ws.Range("K122").Value = 3458.25
$ws.Range("L122").Value = 2720.0001
$ws.Range("M122").Value = -1008.25
$ws.Range("N122").Value = -7620.0001
$ws.Range("H132").Value = 6414376
$ws.Range("I132").Value = 3963.8096
$ws.Range("J132").Value = 33338108
$ws.Range("K132").Value = 11891.4288
$ws.Range("L132").Value = 100014324
$ws.Range("M132").Value = -9361.4288
$ws.Range("N132").Value = -100019384
$ws.Range("H134").Value = 2767.8647
$ws.Range("I134").Value = 2600.3823
$ws.Range("J134").Value = 4666
$ws.Range("K134").Value = 7801.146900000001
$ws.Range("L134").Value = 13998
$ws.Range("M134").Value = -5266.146900000001
$ws.Range("N134").Value = -19068
$ws.Range("H136").Value = 2092.0322
$ws.Range("I136").Value = 1495.5555
$ws.Range("J136").Value = 2917.923
$ws.Range("K136").Value = 4486.666499999999
$ws.Range("L136").Value = 8753.769
$ws.Range("M136").Value = -1936.666499999999
$ws.Range("N136").Value = -13853.769

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 2359.5
$ws.Range("J86").Value = 2379.4443
$ws.Range("L86").Value = 7138.3329
$ws.Range("N86").Value = -9510.332900000001
$ws.Range("H89").Value = 2359.5
$ws.Range("J89").Value = 2379.4443
$ws.Range("L89").Value = 21414.9987
$ws.Range("N89").Value = -33270.9987
$ws.Range("H92").Value = 1710.2222
$ws.Range("J92").Value = 1899.7142
$ws.Range("L92").Value = 5699.142599999999
$ws.Range("N92").Value = -8195.142599999999
$ws.Range("H97").Value = 558052.9
$ws.Range("J97").Value = 3699.3333
$ws.Range("L97").Value = 11097.9999
$ws.Range("N97").Value = -12089.9999
$ws.Range("H98").Value = 588.4375
$ws.Range("J98").Value = 607
$ws.Range("L98").Value = 1821
$ws.Range("N98").Value = -4817

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").Value = ""
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").Value = ""
$ws.Range("H122").Value = 3891.9312
$ws.Range("I122").Value = 2985.7896
$ws.Range("J122").Value = 5613.6
$ws.Range("K122").Value = 8957.3688
$ws.Range("L122").Value = 16840.8
$ws.Range("M122").Value = -6507.3688
$ws.Range("N122").Value = -21740.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2381.9167
$ws.Range("I61").Value = 2343.4285
$ws.Range("K61").Value = 2343.4285
$ws.Range("M61").Value = -2141.4285
$ws.Range("H93").Value = 2248.5
$ws.Range("I93").Value = 2640.2856
$ws.Range("J93").Value = 1334.3334
$ws.Range("K93").Value = 2640.2856
$ws.Range("L93").Value = 1334.3334
$ws.Range("M93").Value = -1392.2856
$ws.Range("N93").Value = -3830.3334
$ws.Range("H113").Value = 2381.9167
$ws.Range("I113").Value = 2343.4285
$ws.Range("K113").Value = 2343.4285
$ws.Range("M113").Value = -173.4285
$ws.Range("H132").Value = 4088.28
$ws.Range("I132").Value = 3798.8
$ws.Range("J132").Value = 4522.5
$ws.Range("K132").Value = 11396.4
$ws.Range("L132").Value = 13567.5
$ws.Range("M132").Value = -8866.400000000001
$ws.Range("N132").Value = -18627.5
$ws.Range("H136").Value = 2634.3635
$ws.Range("I136").Value = 1597.6
$ws.Range("K136").Value = 4792.799999999999
$ws.Range("M136").Value = -2242.799999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 11613.667
$ws.Range("I52").Value = 11613.667
$ws.Range("K52").Value = 11613.667
$ws.Range("M52").Value = -11387.667
$ws.Range("H54").Value = 47034.637
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 47034.637
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 47034.637
$ws.Range("M54").Value = ""
$ws.Range("N54").Value = -48074.637
$ws.Range("H81").Value = 5359.7856
$ws.Range("I81").Value = 4099.273
$ws.Range("J81").Value = 9981.666999999999
$ws.Range("K81").Value = 8198.546
$ws.Range("L81").Value = 19963.334
$ws.Range("M81").Value = -7137.546
$ws.Range("N81").Value = -22085.334
$ws.Range("H84").Value = 5359.7856
$ws.Range("I84").Value = 4099.273
$ws.Range("J84").Value = 9981.666999999999
$ws.Range("K84").Value = 40992.73
$ws.Range("L84").Value = 99816.67
$ws.Range("M84").Value = -35688.73
$ws.Range("N84").Value = -110424.67
$ws.Range("H132").Value = 6132.8237
$ws.Range("J132").Value = 8499.5
$ws.Range("L132").Value = 25498.5
$ws.Range("N132").Value = -30558.5
$ws.Range("H136").Value = 215811.39
$ws.Range("J136").Value = 983073.4
$ws.Range("L136").Value = 2949220.2
$ws.Range("N136").Value = -2954320.2
